# "Fixed and stylized charts"
#
# Underlying data changes on Sheet4 ("Property" / "Percent" summary table):
#   - Row 1 used to hold a lone year value (B1 = 2020); it is replaced with a
#     two-column header "Property" | "Percent".
#   - The remaining rows shift up to follow the new header and the percentage
#     values are re-rounded to 2 decimals (e.g. 51,899047383 -> 51,89).
#   - Sheet4 becomes the active/selected sheet & tab (previously Sheet2 was).
#   - Sheet2's own selection moves off of its data cell onto the last column
#     of its header row.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- Sheet4: rebuild the small summary table -------------------------------
$ws4.Range("B1").Value = "Percent"
$ws4.Range("A1").Value = "Property"

$ws4.Range("A2").Value = "Total of apartments "
$ws4.Range("B2").Value = "51,89"

$ws4.Range("A3").Value = "Total of houses "
$ws4.Range("B3").Value = "46,74"

$ws4.Range("A4").Value = "Total of land"
$ws4.Range("B4").Value = "1,36"

# --- Selection / active-sheet bookkeeping -----------------------------------
# Sheet2 keeps its old selected cell row, but the cursor moves to the last
# (G) column of the header row instead of B2.
$ws2.Range("G2").Select() | Out-Null

# Sheet4 becomes the active tab/sheet, with the cursor parked on the new
# "Percent" header cell.
$ws4.Activate() | Out-Null
$ws4.Range("B1").Select() | Out-Null
